$wb = $excel.ActiveWorkbook

# "Ready for handoff" -> "In Translation" (shared string used on Overview!E2:F2,
# zh-cn!C2 and de-de!C2 - setting the cell Value updates the shared text
# everywhere it is referenced).
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# Narrow the "Status" columns (Overview E & F, zh-cn/de-de column C) to match
# the report's new autofit width.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
